$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.235146045684814
$ws.Range("B1").Value = 2.179468154907227
$ws.Range("C1").Value = 4.203808784484863
$ws.Range("D1").Value = 3.040886878967285
$ws.Range("E1").Value = 1.068594217300415
